$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17, shifting existing rows 17-25 down to 18-26
$ws.Rows.Item(17).Insert(1)

# Populate the new row 17 with data (a duplicate-like entry, same as old row 17/new row 18
# except for the date and volume)
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C17").Value = "Los Lagos"
$ws.Range("D17").Value = 44895
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = "Chirimoya"
$ws.Range("K17").Value = "Cultivar IV Región"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 22000
$ws.Range("O17").Value = 22500
$ws.Range("P17").Value = 22250
$ws.Range("Q17").Value = "$/bandeja 8 kilos"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 2781
$ws.Range("T17").Value = 8
